$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A for brand-new rows 243, 244, 245: set sequence index value and copy style from an existing A-column cell
$ws.Range("A243").Value = 241
$ws.Range("A244").Value = 242
$ws.Range("A245").Value = 243
$ws.Range("A243").Style = $ws.Range("A242").Style
$ws.Range("A244").Style = $ws.Range("A242").Style
$ws.Range("A245").Style = $ws.Range("A242").Style

# Row 179
$ws.Cells.Item(179, 2).Value = 4
$ws.Cells.Item(179, 3).Value = 2000
$ws.Cells.Item(179, 4).Value = 1179
$ws.Cells.Item(179, 5).Value = $true
$ws.Cells.Item(179, 6).Value = $false
$ws.Cells.Item(179, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.44198, 0.94572]))'
$ws.Cells.Item(179, 8).Value = 'testModels/refined_db/4/m410/m410.off'
$ws.Cells.Item(179, 9).Value = 0.4426165134704702
$ws.Cells.Item(179, 10).Value = 0.36472700513595
$ws.Cells.Item(179, 11).Value = $false
$ws.Cells.Item(179, 12).Value = $false

# Row 180
$ws.Cells.Item(180, 2).Value = 4
$ws.Cells.Item(180, 3).Value = 2000
$ws.Cells.Item(180, 4).Value = 1009
$ws.Cells.Item(180, 5).Value = $true
$ws.Cells.Item(180, 6).Value = $false
$ws.Cells.Item(180, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.70357, 0.45049]))'
$ws.Cells.Item(180, 8).Value = 'testModels/refined_db/4/m411/m411.off'
$ws.Cells.Item(180, 9).Value = 0.356823088341781
$ws.Cells.Item(180, 10).Value = 0.2742882714080499
$ws.Cells.Item(180, 11).Value = $false
$ws.Cells.Item(180, 12).Value = $false

# Row 181
$ws.Cells.Item(181, 2).Value = 4
$ws.Cells.Item(181, 3).Value = 2000
$ws.Cells.Item(181, 4).Value = 1472
$ws.Cells.Item(181, 5).Value = $true
$ws.Cells.Item(181, 6).Value = $false
$ws.Cells.Item(181, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.12811, 0.51158]))'
$ws.Cells.Item(181, 8).Value = 'testModels/refined_db/4/m412/m412.off'
$ws.Cells.Item(181, 9).Value = 0.02426980825759749
$ws.Cells.Item(181, 10).Value = 0.04766319038249999
$ws.Cells.Item(181, 11).Value = $false
$ws.Cells.Item(181, 12).Value = $true

# Row 182
$ws.Cells.Item(182, 2).Value = 5
$ws.Cells.Item(182, 3).Value = 1996
$ws.Cells.Item(182, 4).Value = 1000
$ws.Cells.Item(182, 5).Value = $true
$ws.Cells.Item(182, 6).Value = $false
$ws.Cells.Item(182, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.95964, 0.89686, 0.975  ]))'
$ws.Cells.Item(182, 8).Value = 'testModels/refined_db/5/m500/m500.off'
$ws.Cells.Item(182, 9).Value = 0.0006964191987781477
$ws.Cells.Item(182, 10).Value = 0.774133244696
$ws.Cells.Item(182, 11).Value = $false
$ws.Cells.Item(182, 12).Value = $false

# Row 183
$ws.Cells.Item(183, 2).Value = 5
$ws.Cells.Item(183, 3).Value = 2000
$ws.Cells.Item(183, 4).Value = 1072
$ws.Cells.Item(183, 5).Value = $true
$ws.Cells.Item(183, 6).Value = $false
$ws.Cells.Item(183, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.67754, 0.975  , 0.67754]))'
$ws.Cells.Item(183, 8).Value = 'testModels/refined_db/5/m501/m501.off'
$ws.Cells.Item(183, 9).Value = 0.0001625845650226333
$ws.Cells.Item(183, 10).Value = 0.4045155493718
$ws.Cells.Item(183, 11).Value = $false
$ws.Cells.Item(183, 12).Value = $false

# Row 184
$ws.Cells.Item(184, 2).Value = 5
$ws.Cells.Item(184, 3).Value = 1996
$ws.Cells.Item(184, 4).Value = 1000
$ws.Cells.Item(184, 5).Value = $true
$ws.Cells.Item(184, 6).Value = $false
$ws.Cells.Item(184, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.67501, 0.975  , 0.67501]))'
$ws.Cells.Item(184, 8).Value = 'testModels/refined_db/5/m502/m502.off'
$ws.Cells.Item(184, 9).Value = 0.002764792734748884
$ws.Cells.Item(184, 10).Value = 0.4013836450465499
$ws.Cells.Item(184, 11).Value = $false
$ws.Cells.Item(184, 12).Value = $false

# Row 185
$ws.Cells.Item(185, 2).Value = 5
$ws.Cells.Item(185, 3).Value = 1952
$ws.Cells.Item(185, 4).Value = 1000
$ws.Cells.Item(185, 5).Value = $true
$ws.Cells.Item(185, 6).Value = $false
$ws.Cells.Item(185, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.8853, 0.975 , 0.5708]))'
$ws.Cells.Item(185, 8).Value = 'testModels/refined_db/5/m503/m503.off'
$ws.Cells.Item(185, 9).Value = 0.002016172263500733
$ws.Cells.Item(185, 10).Value = 0.44607162209075
$ws.Cells.Item(185, 11).Value = $false
$ws.Cells.Item(185, 12).Value = $false

# Row 186
$ws.Cells.Item(186, 2).Value = 5
$ws.Cells.Item(186, 3).Value = 2000
$ws.Cells.Item(186, 4).Value = 1000
$ws.Cells.Item(186, 5).Value = $true
$ws.Cells.Item(186, 6).Value = $false
$ws.Cells.Item(186, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.73697, 0.77026]))'
$ws.Cells.Item(186, 8).Value = 'testModels/refined_db/5/m504/m504.off'
$ws.Cells.Item(186, 9).Value = 0.000021555741850681358506196597
$ws.Cells.Item(186, 10).Value = 0.50407117646765
$ws.Cells.Item(186, 11).Value = $false
$ws.Cells.Item(186, 12).Value = $false

# Row 187
$ws.Cells.Item(187, 2).Value = 5
$ws.Cells.Item(187, 3).Value = 2000
$ws.Cells.Item(187, 4).Value = 1006
$ws.Cells.Item(187, 5).Value = $true
$ws.Cells.Item(187, 6).Value = $false
$ws.Cells.Item(187, 7).Value = '(TrackedArray([0.02521, 0.02544, 0.025  ]), TrackedArray([0.97299, 0.65558, 0.92819]))'
$ws.Cells.Item(187, 8).Value = 'testModels/refined_db/5/m505/m505.off'
$ws.Cells.Item(187, 9).Value = 0.00526397823374702
$ws.Cells.Item(187, 10).Value = 0.5394184393397173
$ws.Cells.Item(187, 11).Value = $false
$ws.Cells.Item(187, 12).Value = $false

# Row 188
$ws.Cells.Item(188, 2).Value = 5
$ws.Cells.Item(188, 3).Value = 2000
$ws.Cells.Item(188, 4).Value = 1008
$ws.Cells.Item(188, 5).Value = $true
$ws.Cells.Item(188, 6).Value = $false
$ws.Cells.Item(188, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.97449, 0.38986, 0.68268]))'
$ws.Cells.Item(188, 8).Value = 'testModels/refined_db/5/m506/m506.off'
$ws.Cells.Item(188, 9).Value = 0.0003049739221834936
$ws.Cells.Item(188, 10).Value = 0.2278390717677797
$ws.Cells.Item(188, 11).Value = $false
$ws.Cells.Item(188, 12).Value = $false

# Row 189
$ws.Cells.Item(189, 2).Value = 5
$ws.Cells.Item(189, 3).Value = 2000
$ws.Cells.Item(189, 4).Value = 1006
$ws.Cells.Item(189, 5).Value = $true
$ws.Cells.Item(189, 6).Value = $false
$ws.Cells.Item(189, 7).Value = '(TrackedArray([0.025  , 0.02504, 0.025  ]), TrackedArray([0.975  , 0.69504, 0.67724]))'
$ws.Cells.Item(189, 8).Value = 'testModels/refined_db/5/m507/m507.off'
$ws.Cells.Item(189, 9).Value = 0.0004915109394796224
$ws.Cells.Item(189, 10).Value = 0.415150557238935
$ws.Cells.Item(189, 11).Value = $false
$ws.Cells.Item(189, 12).Value = $false

# Row 190
$ws.Cells.Item(190, 2).Value = 5
$ws.Cells.Item(190, 3).Value = 1988
$ws.Cells.Item(190, 4).Value = 998
$ws.Cells.Item(190, 5).Value = $true
$ws.Cells.Item(190, 6).Value = $false
$ws.Cells.Item(190, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.80951, 0.62492]))'
$ws.Cells.Item(190, 8).Value = 'testModels/refined_db/5/m508/m508.off'
$ws.Cells.Item(190, 9).Value = 0.000025932996979407279576753872
$ws.Cells.Item(190, 10).Value = 0.4471125239474999
$ws.Cells.Item(190, 11).Value = $false
$ws.Cells.Item(190, 12).Value = $false

# Row 191
$ws.Cells.Item(191, 2).Value = 5
$ws.Cells.Item(191, 3).Value = 1999
$ws.Cells.Item(191, 4).Value = 1002
$ws.Cells.Item(191, 5).Value = $true
$ws.Cells.Item(191, 6).Value = $false
$ws.Cells.Item(191, 7).Value = '(TrackedArray([0.0254, 0.025 , 0.025 ]), TrackedArray([0.97303, 0.88927, 0.71394]))'
$ws.Cells.Item(191, 8).Value = 'testModels/refined_db/5/m509/m509.off'
$ws.Cells.Item(191, 9).Value = 0.001541805271098911
$ws.Cells.Item(191, 10).Value = 0.5642497394150771
$ws.Cells.Item(191, 11).Value = $false
$ws.Cells.Item(191, 12).Value = $false

# Row 192
$ws.Cells.Item(192, 2).Value = 5
$ws.Cells.Item(192, 3).Value = 1999
$ws.Cells.Item(192, 4).Value = 1035
$ws.Cells.Item(192, 5).Value = $true
$ws.Cells.Item(192, 6).Value = $false
$ws.Cells.Item(192, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.60961, 0.97479, 0.52201]))'
$ws.Cells.Item(192, 8).Value = 'testModels/refined_db/5/m510/m510.off'
$ws.Cells.Item(192, 9).Value = 0.0004058991979104314
$ws.Cells.Item(192, 10).Value = 0.2759672917384885
$ws.Cells.Item(192, 11).Value = $false
$ws.Cells.Item(192, 12).Value = $false

# Row 193
$ws.Cells.Item(193, 2).Value = 5
$ws.Cells.Item(193, 3).Value = 2202
$ws.Cells.Item(193, 4).Value = 1000
$ws.Cells.Item(193, 5).Value = $true
$ws.Cells.Item(193, 6).Value = $false
$ws.Cells.Item(193, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.9522 , 0.975  , 0.92067]))'
$ws.Cells.Item(193, 8).Value = 'testModels/refined_db/5/m511/m511.off'
$ws.Cells.Item(193, 9).Value = 0.02789926138668109
$ws.Cells.Item(193, 10).Value = 0.7889463070987
$ws.Cells.Item(193, 11).Value = $false
$ws.Cells.Item(193, 12).Value = $false

# Row 194
$ws.Cells.Item(194, 2).Value = 5
$ws.Cells.Item(194, 3).Value = 2000
$ws.Cells.Item(194, 4).Value = 1000
$ws.Cells.Item(194, 5).Value = $true
$ws.Cells.Item(194, 6).Value = $false
$ws.Cells.Item(194, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.79759, 0.975  , 0.79759]))'
$ws.Cells.Item(194, 8).Value = 'testModels/refined_db/5/m512/m512.off'
$ws.Cells.Item(194, 9).Value = 0.003047668860747743
$ws.Cells.Item(194, 10).Value = 0.56704613894055
$ws.Cells.Item(194, 11).Value = $false
$ws.Cells.Item(194, 12).Value = $false

# Row 195
$ws.Cells.Item(195, 2).Value = 6
$ws.Cells.Item(195, 3).Value = 2000
$ws.Cells.Item(195, 4).Value = 1042
$ws.Cells.Item(195, 5).Value = $true
$ws.Cells.Item(195, 6).Value = $false
$ws.Cells.Item(195, 7).Value = '(TrackedArray([0.02693, 0.02586, 0.02617]), TrackedArray([0.48928, 0.76584, 0.97387]))'
$ws.Cells.Item(195, 8).Value = 'testModels/refined_db/6/m600/m600.off'
$ws.Cells.Item(195, 9).Value = 0.002188264251238722
$ws.Cells.Item(195, 10).Value = 0.3242315189462048
$ws.Cells.Item(195, 11).Value = $false
$ws.Cells.Item(195, 12).Value = $false

# Row 196
$ws.Cells.Item(196, 2).Value = 6
$ws.Cells.Item(196, 3).Value = 1999
$ws.Cells.Item(196, 4).Value = 1175
$ws.Cells.Item(196, 5).Value = $true
$ws.Cells.Item(196, 6).Value = $false
$ws.Cells.Item(196, 7).Value = '(TrackedArray([0.02796, 0.02501, 0.02509]), TrackedArray([0.975  , 0.93453, 0.95622]))'
$ws.Cells.Item(196, 8).Value = 'testModels/refined_db/6/m601/m601.off'
$ws.Cells.Item(196, 9).Value = 0.001129485291517668
$ws.Cells.Item(196, 10).Value = 0.8020325620060286
$ws.Cells.Item(196, 11).Value = $false
$ws.Cells.Item(196, 12).Value = $true

# Row 197
$ws.Cells.Item(197, 2).Value = 6
$ws.Cells.Item(197, 3).Value = 2000
$ws.Cells.Item(197, 4).Value = 1008
$ws.Cells.Item(197, 5).Value = $true
$ws.Cells.Item(197, 6).Value = $false
$ws.Cells.Item(197, 7).Value = '(TrackedArray([0.02509, 0.025  , 0.025  ]), TrackedArray([0.4276 , 0.975  , 0.59325]))'
$ws.Cells.Item(197, 8).Value = 'testModels/refined_db/6/m602/m602.off'
$ws.Cells.Item(197, 9).Value = 0.0007723776362440639
$ws.Cells.Item(197, 10).Value = 0.217290235051875
$ws.Cells.Item(197, 11).Value = $false
$ws.Cells.Item(197, 12).Value = $false

# Row 198
$ws.Cells.Item(198, 2).Value = 6
$ws.Cells.Item(198, 3).Value = 1996
$ws.Cells.Item(198, 4).Value = 1000
$ws.Cells.Item(198, 5).Value = $true
$ws.Cells.Item(198, 6).Value = $false
$ws.Cells.Item(198, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.58332, 0.84655, 0.975  ]))'
$ws.Cells.Item(198, 8).Value = 'testModels/refined_db/6/m603/m603.off'
$ws.Cells.Item(198, 9).Value = 0.000745184319603115
$ws.Cells.Item(198, 10).Value = 0.4357544367391
$ws.Cells.Item(198, 11).Value = $false
$ws.Cells.Item(198, 12).Value = $false

# Row 199
$ws.Cells.Item(199, 2).Value = 6
$ws.Cells.Item(199, 3).Value = 1999
$ws.Cells.Item(199, 4).Value = 1149
$ws.Cells.Item(199, 5).Value = $true
$ws.Cells.Item(199, 6).Value = $false
$ws.Cells.Item(199, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.69906, 0.975  , 0.15606]))'
$ws.Cells.Item(199, 8).Value = 'testModels/refined_db/6/m604/m604.off'
$ws.Cells.Item(199, 9).Value = 0.02392622713249805
$ws.Cells.Item(199, 10).Value = 0.083922626992
$ws.Cells.Item(199, 11).Value = $false
$ws.Cells.Item(199, 12).Value = $true

# Row 200
$ws.Cells.Item(200, 2).Value = 6
$ws.Cells.Item(200, 3).Value = 2238
$ws.Cells.Item(200, 4).Value = 834
$ws.Cells.Item(200, 5).Value = $true
$ws.Cells.Item(200, 6).Value = $false
$ws.Cells.Item(200, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.50208, 0.975  , 0.83687]))'
$ws.Cells.Item(200, 8).Value = 'testModels/refined_db/6/m605/m605.off'
$ws.Cells.Item(200, 9).Value = 0.002287657914562828
$ws.Cells.Item(200, 10).Value = 0.36795691912105
$ws.Cells.Item(200, 11).Value = $true
$ws.Cells.Item(200, 12).Value = $false

# Row 201
$ws.Cells.Item(201, 2).Value = 6
$ws.Cells.Item(201, 3).Value = 1971
$ws.Cells.Item(201, 4).Value = 1000
$ws.Cells.Item(201, 5).Value = $true
$ws.Cells.Item(201, 6).Value = $false
$ws.Cells.Item(201, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.83864, 0.975  , 0.45899]))'
$ws.Cells.Item(201, 8).Value = 'testModels/refined_db/6/m606/m606.off'
$ws.Cells.Item(201, 9).Value = 0.01730747197990135
$ws.Cells.Item(201, 10).Value = 0.3354520743954
$ws.Cells.Item(201, 11).Value = $false
$ws.Cells.Item(201, 12).Value = $false

# Row 202
$ws.Cells.Item(202, 2).Value = 6
$ws.Cells.Item(202, 3).Value = 1709
$ws.Cells.Item(202, 4).Value = 901
$ws.Cells.Item(202, 5).Value = $true
$ws.Cells.Item(202, 6).Value = $false
$ws.Cells.Item(202, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.34539, 0.975  , 0.84361]))'
$ws.Cells.Item(202, 8).Value = 'testModels/refined_db/6/m607/m607.off'
$ws.Cells.Item(202, 9).Value = 0.000015957809160461080801507103
$ws.Cells.Item(202, 10).Value = 0.2491597883831999
$ws.Cells.Item(202, 11).Value = $false
$ws.Cells.Item(202, 12).Value = $false

# Row 203
$ws.Cells.Item(203, 2).Value = 6
$ws.Cells.Item(203, 3).Value = 1916
$ws.Cells.Item(203, 4).Value = 1000
$ws.Cells.Item(203, 5).Value = $true
$ws.Cells.Item(203, 6).Value = $false
$ws.Cells.Item(203, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.75998, 0.90819, 0.975  ]))'
$ws.Cells.Item(203, 8).Value = 'testModels/refined_db/6/m608/m608.off'
$ws.Cells.Item(203, 9).Value = 0.000086322359758863559755538708
$ws.Cells.Item(203, 10).Value = 0.6166739930119999
$ws.Cells.Item(203, 11).Value = $false
$ws.Cells.Item(203, 12).Value = $false

# Row 204
$ws.Cells.Item(204, 2).Value = 6
$ws.Cells.Item(204, 3).Value = 2000
$ws.Cells.Item(204, 4).Value = 1020
$ws.Cells.Item(204, 5).Value = $true
$ws.Cells.Item(204, 6).Value = $false
$ws.Cells.Item(204, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.21256, 0.72068, 0.975  ]))'
$ws.Cells.Item(204, 8).Value = 'testModels/refined_db/6/m609/m609.off'
$ws.Cells.Item(204, 9).Value = 0.0003560311298754342
$ws.Cells.Item(204, 10).Value = 0.1239595846008
$ws.Cells.Item(204, 11).Value = $false
$ws.Cells.Item(204, 12).Value = $false

# Row 205
$ws.Cells.Item(205, 2).Value = 6
$ws.Cells.Item(205, 3).Value = 2000
$ws.Cells.Item(205, 4).Value = 1016
$ws.Cells.Item(205, 5).Value = $true
$ws.Cells.Item(205, 6).Value = $false
$ws.Cells.Item(205, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.6368 , 0.55179]))'
$ws.Cells.Item(205, 8).Value = 'testModels/refined_db/6/m610/m610.off'
$ws.Cells.Item(205, 9).Value = 0.000013313841439356409847641349
$ws.Cells.Item(205, 10).Value = 0.3061758177928125
$ws.Cells.Item(205, 11).Value = $false
$ws.Cells.Item(205, 12).Value = $false

# Row 206
$ws.Cells.Item(206, 2).Value = 6
$ws.Cells.Item(206, 3).Value = 2026
$ws.Cells.Item(206, 4).Value = 778
$ws.Cells.Item(206, 5).Value = $true
$ws.Cells.Item(206, 6).Value = $false
$ws.Cells.Item(206, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.84401, 0.975  , 0.45893]))'
$ws.Cells.Item(206, 8).Value = 'testModels/refined_db/6/m611/m611.off'
$ws.Cells.Item(206, 9).Value = 0.006340770550030407
$ws.Cells.Item(206, 10).Value = 0.337626976131
$ws.Cells.Item(206, 11).Value = $true
$ws.Cells.Item(206, 12).Value = $false

# Row 207
$ws.Cells.Item(207, 2).Value = 6
$ws.Cells.Item(207, 3).Value = 2232
$ws.Cells.Item(207, 4).Value = 948
$ws.Cells.Item(207, 5).Value = $true
$ws.Cells.Item(207, 6).Value = $false
$ws.Cells.Item(207, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.75993, 0.90811, 0.975  ]))'
$ws.Cells.Item(207, 8).Value = 'testModels/refined_db/6/m612/m612.off'
$ws.Cells.Item(207, 9).Value = 0.0001209419976838267
$ws.Cells.Item(207, 10).Value = 0.6165689174602
$ws.Cells.Item(207, 11).Value = $false
$ws.Cells.Item(207, 12).Value = $false

# Row 208
$ws.Cells.Item(208, 2).Value = 7
$ws.Cells.Item(208, 3).Value = 1904
$ws.Cells.Item(208, 4).Value = 1000
$ws.Cells.Item(208, 5).Value = $true
$ws.Cells.Item(208, 6).Value = $false
$ws.Cells.Item(208, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.41837, 0.975  , 0.09773]))'
$ws.Cells.Item(208, 8).Value = 'testModels/refined_db/7/m700/m700.off'
$ws.Cells.Item(208, 9).Value = 0.001634965579034274
$ws.Cells.Item(208, 10).Value = 0.027180761887565
$ws.Cells.Item(208, 11).Value = $false
$ws.Cells.Item(208, 12).Value = $false

# Row 209
$ws.Cells.Item(209, 2).Value = 7
$ws.Cells.Item(209, 3).Value = 1900
$ws.Cells.Item(209, 4).Value = 960
$ws.Cells.Item(209, 5).Value = $true
$ws.Cells.Item(209, 6).Value = $false
$ws.Cells.Item(209, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.04405, 0.30701, 0.975  ]))'
$ws.Cells.Item(209, 8).Value = 'testModels/refined_db/7/m701/m701.off'
$ws.Cells.Item(209, 9).Value = 0.0005984744284004732
$ws.Cells.Item(209, 10).Value = 0.005102765082699999
$ws.Cells.Item(209, 11).Value = $false
$ws.Cells.Item(209, 12).Value = $false

# Row 210
$ws.Cells.Item(210, 2).Value = 7
$ws.Cells.Item(210, 3).Value = 2000
$ws.Cells.Item(210, 4).Value = 1005
$ws.Cells.Item(210, 5).Value = $true
$ws.Cells.Item(210, 6).Value = $false
$ws.Cells.Item(210, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.12358, 0.12358, 0.975  ]))'
$ws.Cells.Item(210, 8).Value = 'testModels/refined_db/7/m702/m702.off'
$ws.Cells.Item(210, 9).Value = 0.0002847790346095218
$ws.Cells.Item(210, 10).Value = 0.0092317409798
$ws.Cells.Item(210, 11).Value = $false
$ws.Cells.Item(210, 12).Value = $false

# Row 211
$ws.Cells.Item(211, 2).Value = 7
$ws.Cells.Item(211, 3).Value = 1924
$ws.Cells.Item(211, 4).Value = 1000
$ws.Cells.Item(211, 5).Value = $true
$ws.Cells.Item(211, 6).Value = $false
$ws.Cells.Item(211, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.05076, 0.12084, 0.975  ]))'
$ws.Cells.Item(211, 8).Value = 'testModels/refined_db/7/m703/m703.off'
$ws.Cells.Item(211, 9).Value = 0.0001905921691495498
$ws.Cells.Item(211, 10).Value = 0.002345430056895
$ws.Cells.Item(211, 11).Value = $false
$ws.Cells.Item(211, 12).Value = $false

# Row 212
$ws.Cells.Item(212, 2).Value = 7
$ws.Cells.Item(212, 3).Value = 1996
$ws.Cells.Item(212, 4).Value = 1000
$ws.Cells.Item(212, 5).Value = $true
$ws.Cells.Item(212, 6).Value = $false
$ws.Cells.Item(212, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.05493, 0.08613, 0.975  ]))'
$ws.Cells.Item(212, 8).Value = 'testModels/refined_db/7/m704/m704.off'
$ws.Cells.Item(212, 9).Value = 0.000090252913008419981873847737
$ws.Cells.Item(212, 10).Value = 0.001738359126875
$ws.Cells.Item(212, 11).Value = $false
$ws.Cells.Item(212, 12).Value = $false

# Row 213
$ws.Cells.Item(213, 2).Value = 7
$ws.Cells.Item(213, 3).Value = 2000
$ws.Cells.Item(213, 4).Value = 1028
$ws.Cells.Item(213, 5).Value = $true
$ws.Cells.Item(213, 6).Value = $false
$ws.Cells.Item(213, 7).Value = '(TrackedArray([0.02655, 0.025  , 0.02525]), TrackedArray([0.46452, 0.975  , 0.09452]))'
$ws.Cells.Item(213, 8).Value = 'testModels/refined_db/7/m705/m705.off'
$ws.Cells.Item(213, 9).Value = 0.0005383789458012969
$ws.Cells.Item(213, 10).Value = 0.02882194424788575
$ws.Cells.Item(213, 11).Value = $false
$ws.Cells.Item(213, 12).Value = $false

# Row 214
$ws.Cells.Item(214, 2).Value = 7
$ws.Cells.Item(214, 3).Value = 1999
$ws.Cells.Item(214, 4).Value = 802
$ws.Cells.Item(214, 5).Value = $true
$ws.Cells.Item(214, 6).Value = $false
$ws.Cells.Item(214, 7).Value = '(TrackedArray([0.025 , 0.025 , 0.0252]), TrackedArray([0.28441, 0.975  , 0.08598]))'
$ws.Cells.Item(214, 8).Value = 'testModels/refined_db/7/m706/m706.off'
$ws.Cells.Item(214, 9).Value = 0.000843582869486666
$ws.Cells.Item(214, 10).Value = 0.01497704798992
$ws.Cells.Item(214, 11).Value = $true
$ws.Cells.Item(214, 12).Value = $false

# Row 215
$ws.Cells.Item(215, 2).Value = 7
$ws.Cells.Item(215, 3).Value = 1999
$ws.Cells.Item(215, 4).Value = 1005
$ws.Cells.Item(215, 5).Value = $true
$ws.Cells.Item(215, 6).Value = $false
$ws.Cells.Item(215, 7).Value = '(TrackedArray([0.02506, 0.025  , 0.02503]), TrackedArray([0.26975, 0.975  , 0.06549]))'
$ws.Cells.Item(215, 8).Value = 'testModels/refined_db/7/m707/m707.off'
$ws.Cells.Item(215, 9).Value = 0.06494472469930278
$ws.Cells.Item(215, 10).Value = 0.009405908309349748
$ws.Cells.Item(215, 11).Value = $false
$ws.Cells.Item(215, 12).Value = $false

# Row 216
$ws.Cells.Item(216, 2).Value = 7
$ws.Cells.Item(216, 3).Value = 2000
$ws.Cells.Item(216, 4).Value = 1032
$ws.Cells.Item(216, 5).Value = $true
$ws.Cells.Item(216, 6).Value = $false
$ws.Cells.Item(216, 7).Value = '(TrackedArray([0.02512, 0.025  , 0.025  ]), TrackedArray([0.44415, 0.975  , 0.06419]))'
$ws.Cells.Item(216, 8).Value = 'testModels/refined_db/7/m708/m708.off'
$ws.Cells.Item(216, 9).Value = 0.004282833850689033
$ws.Cells.Item(216, 10).Value = 0.01559969183422879
$ws.Cells.Item(216, 11).Value = $false
$ws.Cells.Item(216, 12).Value = $false

# Row 217
$ws.Cells.Item(217, 2).Value = 7
$ws.Cells.Item(217, 3).Value = 2000
$ws.Cells.Item(217, 4).Value = 1014
$ws.Cells.Item(217, 5).Value = $true
$ws.Cells.Item(217, 6).Value = $false
$ws.Cells.Item(217, 7).Value = '(TrackedArray([0.025  , 0.02539, 0.025  ]), TrackedArray([0.22979, 0.07604, 0.97442]))'
$ws.Cells.Item(217, 8).Value = 'testModels/refined_db/7/m709/m709.off'
$ws.Cells.Item(217, 9).Value = 0.0003573684656633618
$ws.Cells.Item(217, 10).Value = 0.009848593029662593
$ws.Cells.Item(217, 11).Value = $false
$ws.Cells.Item(217, 12).Value = $false

# Row 218
$ws.Cells.Item(218, 2).Value = 7
$ws.Cells.Item(218, 3).Value = 2000
$ws.Cells.Item(218, 4).Value = 1048
$ws.Cells.Item(218, 5).Value = $true
$ws.Cells.Item(218, 6).Value = $false
$ws.Cells.Item(218, 7).Value = '(TrackedArray([0.0255 , 0.02814, 0.025  ]), TrackedArray([0.22092, 0.07791, 0.97344]))'
$ws.Cells.Item(218, 8).Value = 'testModels/refined_db/7/m710/m710.off'
$ws.Cells.Item(218, 9).Value = 0.00977703969522959
$ws.Cells.Item(218, 10).Value = 0.009225152976660354
$ws.Cells.Item(218, 11).Value = $false
$ws.Cells.Item(218, 12).Value = $false

# Row 219
$ws.Cells.Item(219, 2).Value = 7
$ws.Cells.Item(219, 3).Value = 1999
$ws.Cells.Item(219, 4).Value = 1152
$ws.Cells.Item(219, 5).Value = $true
$ws.Cells.Item(219, 6).Value = $false
$ws.Cells.Item(219, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.22872, 0.0556 , 0.97108]))'
$ws.Cells.Item(219, 8).Value = 'testModels/refined_db/7/m711/m711.off'
$ws.Cells.Item(219, 9).Value = 0.000018726054740883578913048305
$ws.Cells.Item(219, 10).Value = 0.005898091483909225
$ws.Cells.Item(219, 11).Value = $false
$ws.Cells.Item(219, 12).Value = $true

# Row 220
$ws.Cells.Item(220, 2).Value = 7
$ws.Cells.Item(220, 3).Value = 2000
$ws.Cells.Item(220, 4).Value = 1096
$ws.Cells.Item(220, 5).Value = $true
$ws.Cells.Item(220, 6).Value = $false
$ws.Cells.Item(220, 7).Value = '(TrackedArray([0.025  , 0.025  , 0.02511]), TrackedArray([0.97497, 0.0565 , 0.2089 ]))'
$ws.Cells.Item(220, 8).Value = 'testModels/refined_db/7/m712/m712.off'
$ws.Cells.Item(220, 9).Value = 0.0008682165576076732
$ws.Cells.Item(220, 10).Value = 0.005499642830204251
$ws.Cells.Item(220, 11).Value = $false
$ws.Cells.Item(220, 12).Value = $false

# Row 221
$ws.Cells.Item(221, 2).Value = 8
$ws.Cells.Item(221, 3).Value = 1958
$ws.Cells.Item(221, 4).Value = 1000
$ws.Cells.Item(221, 5).Value = $true
$ws.Cells.Item(221, 6).Value = $false
$ws.Cells.Item(221, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.72778, 0.975  , 0.72778]))'
$ws.Cells.Item(221, 8).Value = 'testModels/refined_db/8/m800/m800.off'
$ws.Cells.Item(221, 9).Value = 0.00229425835594984
$ws.Cells.Item(221, 10).Value = 0.4692020714198
$ws.Cells.Item(221, 11).Value = $false
$ws.Cells.Item(221, 12).Value = $false

# Row 222
$ws.Cells.Item(222, 2).Value = 8
$ws.Cells.Item(222, 3).Value = 2000
$ws.Cells.Item(222, 4).Value = 1054
$ws.Cells.Item(222, 5).Value = $true
$ws.Cells.Item(222, 6).Value = $false
$ws.Cells.Item(222, 7).Value = '(TrackedArray([0.02682, 0.02561, 0.025  ]), TrackedArray([0.46699, 0.975  , 0.54432]))'
$ws.Cells.Item(222, 8).Value = 'testModels/refined_db/8/m801/m801.off'
$ws.Cells.Item(222, 9).Value = 0.02356440389067259
$ws.Cells.Item(222, 10).Value = 0.2170203488597468
$ws.Cells.Item(222, 11).Value = $false
$ws.Cells.Item(222, 12).Value = $false

# Row 223
$ws.Cells.Item(223, 2).Value = 8
$ws.Cells.Item(223, 3).Value = 1850
$ws.Cells.Item(223, 4).Value = 1000
$ws.Cells.Item(223, 5).Value = $true
$ws.Cells.Item(223, 6).Value = $false
$ws.Cells.Item(223, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.68578, 0.975  , 0.69614]))'
$ws.Cells.Item(223, 8).Value = 'testModels/refined_db/8/m802/m802.off'
$ws.Cells.Item(223, 9).Value = 0.000065096969130359473945136861
$ws.Cells.Item(223, 10).Value = 0.4212982889336
$ws.Cells.Item(223, 11).Value = $false
$ws.Cells.Item(223, 12).Value = $false

# Row 224
$ws.Cells.Item(224, 2).Value = 8
$ws.Cells.Item(224, 3).Value = 2000
$ws.Cells.Item(224, 4).Value = 1058
$ws.Cells.Item(224, 5).Value = $true
$ws.Cells.Item(224, 6).Value = $false
$ws.Cells.Item(224, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.5701 , 0.975  , 0.60681]))'
$ws.Cells.Item(224, 8).Value = 'testModels/refined_db/8/m803/m803.off'
$ws.Cells.Item(224, 9).Value = 0.0001037527882258007
$ws.Cells.Item(224, 10).Value = 0.30128957545635
$ws.Cells.Item(224, 11).Value = $false
$ws.Cells.Item(224, 12).Value = $false

# Row 225
$ws.Cells.Item(225, 2).Value = 8
$ws.Cells.Item(225, 3).Value = 1900
$ws.Cells.Item(225, 4).Value = 1000
$ws.Cells.Item(225, 5).Value = $true
$ws.Cells.Item(225, 6).Value = $false
$ws.Cells.Item(225, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.70218, 0.975  , 0.66904]))'
$ws.Cells.Item(225, 8).Value = 'testModels/refined_db/8/m804/m804.off'
$ws.Cells.Item(225, 9).Value = 0.002552374671082288
$ws.Cells.Item(225, 10).Value = 0.4143275475168
$ws.Cells.Item(225, 11).Value = $false
$ws.Cells.Item(225, 12).Value = $false

# Row 226
$ws.Cells.Item(226, 2).Value = 8
$ws.Cells.Item(226, 3).Value = 1992
$ws.Cells.Item(226, 4).Value = 1000
$ws.Cells.Item(226, 5).Value = $true
$ws.Cells.Item(226, 6).Value = $false
$ws.Cells.Item(226, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.72976, 0.975  , 0.72976]))'
$ws.Cells.Item(226, 8).Value = 'testModels/refined_db/8/m805/m805.off'
$ws.Cells.Item(226, 9).Value = 0.0003921532034501675
$ws.Cells.Item(226, 10).Value = 0.4718550028118
$ws.Cells.Item(226, 11).Value = $false
$ws.Cells.Item(226, 12).Value = $false

# Row 227
$ws.Cells.Item(227, 2).Value = 8
$ws.Cells.Item(227, 3).Value = 1968
$ws.Cells.Item(227, 4).Value = 1000
$ws.Cells.Item(227, 5).Value = $true
$ws.Cells.Item(227, 6).Value = $false
$ws.Cells.Item(227, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.35829, 0.975  , 0.41154]))'
$ws.Cells.Item(227, 8).Value = 'testModels/refined_db/8/m806/m806.off'
$ws.Cells.Item(227, 9).Value = 0.000062301265942574980834754717
$ws.Cells.Item(227, 10).Value = 0.12238963667975
$ws.Cells.Item(227, 11).Value = $false
$ws.Cells.Item(227, 12).Value = $false

# Row 228
$ws.Cells.Item(228, 2).Value = 8
$ws.Cells.Item(228, 3).Value = 1976
$ws.Cells.Item(228, 4).Value = 1000
$ws.Cells.Item(228, 5).Value = $true
$ws.Cells.Item(228, 6).Value = $false
$ws.Cells.Item(228, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.405  , 0.975  , 0.40384]))'
$ws.Cells.Item(228, 8).Value = 'testModels/refined_db/8/m807/m807.off'
$ws.Cells.Item(228, 9).Value = 0.001213741904721069
$ws.Cells.Item(228, 10).Value = 0.136760518
$ws.Cells.Item(228, 11).Value = $false
$ws.Cells.Item(228, 12).Value = $false

# Row 229
$ws.Cells.Item(229, 2).Value = 8
$ws.Cells.Item(229, 3).Value = 1660
$ws.Cells.Item(229, 4).Value = 1000
$ws.Cells.Item(229, 5).Value = $true
$ws.Cells.Item(229, 6).Value = $false
$ws.Cells.Item(229, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.55085, 0.975  , 0.62458]))'
$ws.Cells.Item(229, 8).Value = 'testModels/refined_db/8/m808/m808.off'
$ws.Cells.Item(229, 9).Value = 0.0005334409098800224
$ws.Cells.Item(229, 10).Value = 0.29952589509265
$ws.Cells.Item(229, 11).Value = $false
$ws.Cells.Item(229, 12).Value = $false

# Row 230
$ws.Cells.Item(230, 2).Value = 8
$ws.Cells.Item(230, 3).Value = 1940
$ws.Cells.Item(230, 4).Value = 1000
$ws.Cells.Item(230, 5).Value = $true
$ws.Cells.Item(230, 6).Value = $false
$ws.Cells.Item(230, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.24254, 0.975  , 0.22572]))'
$ws.Cells.Item(230, 8).Value = 'testModels/refined_db/8/m809/m809.off'
$ws.Cells.Item(230, 9).Value = 0.001741736583522547
$ws.Cells.Item(230, 10).Value = 0.04148136539820001
$ws.Cells.Item(230, 11).Value = $false
$ws.Cells.Item(230, 12).Value = $false

# Row 231
$ws.Cells.Item(231, 2).Value = 8
$ws.Cells.Item(231, 3).Value = 1944
$ws.Cells.Item(231, 4).Value = 1000
$ws.Cells.Item(231, 5).Value = $true
$ws.Cells.Item(231, 6).Value = $false
$ws.Cells.Item(231, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.38296, 0.975  , 0.3995 ]))'
$ws.Cells.Item(231, 8).Value = 'testModels/refined_db/8/m810/m810.off'
$ws.Cells.Item(231, 9).Value = 0.0005863661447857373
$ws.Cells.Item(231, 10).Value = 0.1273507914384
$ws.Cells.Item(231, 11).Value = $false
$ws.Cells.Item(231, 12).Value = $false

# Row 232
$ws.Cells.Item(232, 2).Value = 8
$ws.Cells.Item(232, 3).Value = 1920
$ws.Cells.Item(232, 4).Value = 964
$ws.Cells.Item(232, 5).Value = $true
$ws.Cells.Item(232, 6).Value = $false
$ws.Cells.Item(232, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.53825, 0.975  , 0.55673]))'
$ws.Cells.Item(232, 8).Value = 'testModels/refined_db/8/m811/m811.off'
$ws.Cells.Item(232, 9).Value = 0.0003918203221181705
$ws.Cells.Item(232, 10).Value = 0.2592673744301999
$ws.Cells.Item(232, 11).Value = $false
$ws.Cells.Item(232, 12).Value = $false

# Row 233
$ws.Cells.Item(233, 2).Value = 8
$ws.Cells.Item(233, 3).Value = 1920
$ws.Cells.Item(233, 4).Value = 976
$ws.Cells.Item(233, 5).Value = $true
$ws.Cells.Item(233, 6).Value = $false
$ws.Cells.Item(233, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.46015, 0.975  , 0.47582]))'
$ws.Cells.Item(233, 8).Value = 'testModels/refined_db/8/m812/m812.off'
$ws.Cells.Item(233, 9).Value = 0.000096377794857304151091882682
$ws.Cells.Item(233, 10).Value = 0.18636813187305
$ws.Cells.Item(233, 11).Value = $false
$ws.Cells.Item(233, 12).Value = $false

# Row 234
$ws.Cells.Item(234, 2).Value = 9
$ws.Cells.Item(234, 3).Value = 1772
$ws.Cells.Item(234, 4).Value = 1000
$ws.Cells.Item(234, 5).Value = $true
$ws.Cells.Item(234, 6).Value = $false
$ws.Cells.Item(234, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.61843, 0.38125]))'
$ws.Cells.Item(234, 8).Value = 'testModels/refined_db/9/m900/m900.off'
$ws.Cells.Item(234, 9).Value = 0.0009292084722142014
$ws.Cells.Item(234, 10).Value = 0.2008387393843999
$ws.Cells.Item(234, 11).Value = $false
$ws.Cells.Item(234, 12).Value = $false

# Row 235
$ws.Cells.Item(235, 2).Value = 9
$ws.Cells.Item(235, 3).Value = 1980
$ws.Cells.Item(235, 4).Value = 1000
$ws.Cells.Item(235, 5).Value = $true
$ws.Cells.Item(235, 6).Value = $false
$ws.Cells.Item(235, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.74465, 0.75526]))'
$ws.Cells.Item(235, 8).Value = 'testModels/refined_db/9/m901/m901.off'
$ws.Cells.Item(235, 9).Value = 0.0008655712038301684
$ws.Cells.Item(235, 10).Value = 0.4992543045531
$ws.Cells.Item(235, 11).Value = $false
$ws.Cells.Item(235, 12).Value = $false

# Row 236
$ws.Cells.Item(236, 2).Value = 9
$ws.Cells.Item(236, 3).Value = 1976
$ws.Cells.Item(236, 4).Value = 1000
$ws.Cells.Item(236, 5).Value = $true
$ws.Cells.Item(236, 6).Value = $false
$ws.Cells.Item(236, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.59618, 0.43685]))'
$ws.Cells.Item(236, 8).Value = 'testModels/refined_db/9/m902/m902.off'
$ws.Cells.Item(236, 9).Value = 0.002001445678988587
$ws.Cells.Item(236, 10).Value = 0.2234783074855499
$ws.Cells.Item(236, 11).Value = $false
$ws.Cells.Item(236, 12).Value = $false

# Row 237
$ws.Cells.Item(237, 2).Value = 9
$ws.Cells.Item(237, 3).Value = 1980
$ws.Cells.Item(237, 4).Value = 1000
$ws.Cells.Item(237, 5).Value = $true
$ws.Cells.Item(237, 6).Value = $false
$ws.Cells.Item(237, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.25027, 0.23617]))'
$ws.Cells.Item(237, 8).Value = 'testModels/refined_db/9/m903/m903.off'
$ws.Cells.Item(237, 9).Value = 0.00108864211424962
$ws.Cells.Item(237, 10).Value = 0.0451911105855
$ws.Cells.Item(237, 11).Value = $false
$ws.Cells.Item(237, 12).Value = $false

# Row 238
$ws.Cells.Item(238, 2).Value = 9
$ws.Cells.Item(238, 3).Value = 1896
$ws.Cells.Item(238, 4).Value = 1000
$ws.Cells.Item(238, 5).Value = $true
$ws.Cells.Item(238, 6).Value = $false
$ws.Cells.Item(238, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.40595, 0.975  , 0.40595]))'
$ws.Cells.Item(238, 8).Value = 'testModels/refined_db/9/m904/m904.off'
$ws.Cells.Item(238, 9).Value = 0.000091530553867136792015393731
$ws.Cells.Item(238, 10).Value = 0.13786748118095
$ws.Cells.Item(238, 11).Value = $false
$ws.Cells.Item(238, 12).Value = $false

# Row 239
$ws.Cells.Item(239, 2).Value = 9
$ws.Cells.Item(239, 3).Value = 1972
$ws.Cells.Item(239, 4).Value = 995
$ws.Cells.Item(239, 5).Value = $true
$ws.Cells.Item(239, 6).Value = $false
$ws.Cells.Item(239, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.84105, 0.975  ]))'
$ws.Cells.Item(239, 8).Value = 'testModels/refined_db/9/m905/m905.off'
$ws.Cells.Item(239, 9).Value = 0.01783165915719444
$ws.Cells.Item(239, 10).Value = 0.7364851249999999
$ws.Cells.Item(239, 11).Value = $false
$ws.Cells.Item(239, 12).Value = $false

# Row 240
$ws.Cells.Item(240, 2).Value = 9
$ws.Cells.Item(240, 3).Value = 1956
$ws.Cells.Item(240, 4).Value = 1000
$ws.Cells.Item(240, 5).Value = $true
$ws.Cells.Item(240, 6).Value = $false
$ws.Cells.Item(240, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.73251, 0.35051, 0.975  ]))'
$ws.Cells.Item(240, 8).Value = 'testModels/refined_db/9/m906/m906.off'
$ws.Cells.Item(240, 9).Value = 0.000093694654164776439825811083
$ws.Cells.Item(240, 10).Value = 0.2187877917004
$ws.Cells.Item(240, 11).Value = $false
$ws.Cells.Item(240, 12).Value = $false

# Row 241
$ws.Cells.Item(241, 2).Value = 9
$ws.Cells.Item(241, 3).Value = 1898
$ws.Cells.Item(241, 4).Value = 1000
$ws.Cells.Item(241, 5).Value = $true
$ws.Cells.Item(241, 6).Value = $false
$ws.Cells.Item(241, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.46091, 0.975  , 0.46112]))'
$ws.Cells.Item(241, 8).Value = 'testModels/refined_db/9/m907/m907.off'
$ws.Cells.Item(241, 9).Value = 0.000000090918065913617563478997
$ws.Cells.Item(241, 10).Value = 0.1806036165228
$ws.Cells.Item(241, 11).Value = $false
$ws.Cells.Item(241, 12).Value = $false

# Row 242
$ws.Cells.Item(242, 2).Value = 9
$ws.Cells.Item(242, 3).Value = 1772
$ws.Cells.Item(242, 4).Value = 1000
$ws.Cells.Item(242, 5).Value = $true
$ws.Cells.Item(242, 6).Value = $false
$ws.Cells.Item(242, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.38093, 0.5    ]))'
$ws.Cells.Item(242, 8).Value = 'testModels/refined_db/9/m908/m908.off'
$ws.Cells.Item(242, 9).Value = 0.0006813556252919533
$ws.Cells.Item(242, 10).Value = 0.1606127350113
$ws.Cells.Item(242, 11).Value = $false
$ws.Cells.Item(242, 12).Value = $false

# Row 243
$ws.Cells.Item(243, 2).Value = 9
$ws.Cells.Item(243, 3).Value = 2000
$ws.Cells.Item(243, 4).Value = 1010
$ws.Cells.Item(243, 5).Value = $true
$ws.Cells.Item(243, 6).Value = $false
$ws.Cells.Item(243, 7).Value = '(TrackedArray([0.025  , 0.02623, 0.025  ]), TrackedArray([0.70357, 0.44378, 0.975  ]))'
$ws.Cells.Item(243, 8).Value = 'testModels/refined_db/9/m909/m909.off'
$ws.Cells.Item(243, 9).Value = 0.1286547542692651
$ws.Cells.Item(243, 10).Value = 0.2691700682120299
$ws.Cells.Item(243, 11).Value = $false
$ws.Cells.Item(243, 12).Value = $false

# Row 244
$ws.Cells.Item(244, 2).Value = 9
$ws.Cells.Item(244, 3).Value = 1974
$ws.Cells.Item(244, 4).Value = 1000
$ws.Cells.Item(244, 5).Value = $true
$ws.Cells.Item(244, 6).Value = $false
$ws.Cells.Item(244, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975  , 0.47863, 0.5    ]))'
$ws.Cells.Item(244, 8).Value = 'testModels/refined_db/9/m911/m911.off'
$ws.Cells.Item(244, 9).Value = 0.004597958658645679
$ws.Cells.Item(244, 10).Value = 0.20469828125
$ws.Cells.Item(244, 11).Value = $false
$ws.Cells.Item(244, 12).Value = $false

# Row 245
$ws.Cells.Item(245, 2).Value = 9
$ws.Cells.Item(245, 3).Value = 1988
$ws.Cells.Item(245, 4).Value = 1000
$ws.Cells.Item(245, 5).Value = $true
$ws.Cells.Item(245, 6).Value = $false
$ws.Cells.Item(245, 7).Value = '(TrackedArray([0.025, 0.025, 0.025]), TrackedArray([0.975, 0.5  , 0.5  ]))'
$ws.Cells.Item(245, 8).Value = 'testModels/refined_db/9/m912/m912.off'
$ws.Cells.Item(245, 9).Value = 0.001119938230630489
$ws.Cells.Item(245, 10).Value = 0.21434375
$ws.Cells.Item(245, 11).Value = $false
$ws.Cells.Item(245, 12).Value = $false
